$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -12.0879
$ws.Range("A3").Value = -21.92459999999999
$ws.Range("D3").Value = -7.257799999999994
$ws.Range("D12").Value = -7.259900000000003
$ws.Range("A14").Value = -21.90619999999999
$ws.Range("A21").Value = -19.96359999999998
$ws.Range("A23").Value = -20.34879999999998
$ws.Range("D24").Value = -7.521699999999998
$ws.Range("A25").Value = -22.01259999999999
$ws.Range("C25").Value = -12.91119999999999
$ws.Range("D25").Value = -8.658799999999996
$ws.Range("A26").Value = -21.03569999999996
$ws.Range("C27").Value = -12.8204
$ws.Range("A29").Value = -21.05289999999998
$ws.Range("C31").Value = -12.8017
$ws.Range("C39").Value = -12.63310000000001
$ws.Range("C48").Value = -11.37889999999999
$ws.Range("D50").Value = -8.142300000000002
$ws.Range("C51").Value = -11.8758
$ws.Range("C52").Value = -11.1242
$ws.Range("A53").Value = -22.30500000000001
$ws.Range("D53").Value = -6.0726
$ws.Range("C55").Value = -13.57539999999999
$ws.Range("C56").Value = -12.1467
$ws.Range("A57").Value = -21.91980000000002
$ws.Range("C57").Value = -12.69369999999999
$ws.Range("D57").Value = -8.659799999999997
$ws.Range("A59").Value = -22.34729999999999
$ws.Range("D61").Value = -7.892899999999996
$ws.Range("D63").Value = -8.123699999999999
$ws.Range("A69").Value = -21.58000000000001
$ws.Range("D70").Value = -7.340399999999996
$ws.Range("C73").Value = -12.87190000000001
$ws.Range("A79").Value = -20.46250000000002
$ws.Range("A83").Value = -21.94749999999999
$ws.Range("D86").Value = -7.612899999999997
$ws.Range("C89").Value = -10.3603
$ws.Range("C90").Value = -12.2328
$ws.Range("A91").Value = -21.36370000000003
$ws.Range("C92").Value = -10.7008
$ws.Range("A93").Value = -20.972
$ws.Range("D98").Value = -8.388200000000005
$ws.Range("D100").Value = -8.173100000000005
$ws.Range("D102").Value = -7.858699999999996
